# Yahoo Finance workbook update:
# Adds several new "scratch" cells (I5, L5, L9, C17, E18 formula, E19, G19,
# G21, B22, C22, F25) to Sheet1, extending the used range from A1:G12 to
# A1:L25, and moves the active selection to G15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scratch values scattered around the sheet (beyond the original
# A1:G12 stock-quote table).
$ws.Range("I5").Value  = "eqww"
$ws.Range("L5").Value  = "asdqwe"
$ws.Range("L9").Value  = "XD"

$ws.Range("C17").Value = "qweqwe"

# A little arithmetic formula plus a label pointing at it.
$ws.Range("E18").Formula = "=B4+C4"
$ws.Range("E19").Value  = "^ this is a formula"
$ws.Range("G19").Value  = "XDDD"

$ws.Range("G21").Value  = "qeqweq"

$ws.Range("B22").Value  = "eqw"
$ws.Range("C22").Value  = "qweqwewq"

$ws.Range("F25").Value  = "eqweqw"

# Leave the active cell on G15, matching the saved selection state.
$ws.Range("G15").Select()
